# ---------------------------------------------------------------------------
# Edit script: updates MasterTestData.xlsx for "scenario 1"
#
#  Sheet "channel_management" (Worksheets.Item(1))
#    - the small 2-column / 6-row "Scenario" lookup table (Table2) is
#      replaced by a wide 9-column / 1-data-row table describing the TV
#      channel default scenario (ScenarioID, ChannelNameTV, ...).
#
#  Sheet "comparison_management" (Worksheets.Item(2))
#    - the lookup table's first column is renamed from "Scenario" to
#      "ScenarioID" (text of the other rows is unchanged).
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # channel_management
$ws2 = $wb.Worksheets.Item(2)   # comparison_management

# ---------------------------------------------------------------------------
# 1. channel_management - remove the old 5 extra data rows (rows 3..6),
#    keeping just the header row and a single data row.
# ---------------------------------------------------------------------------
$ws1.Range("A3:B6").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Resize the worksheet's table (Table2) to the new 9 column x 2 row
#    range *before* writing the new header text, so that each ListColumn
#    lines up with the header cell that will carry its name.
# ---------------------------------------------------------------------------
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:I2"))

# ---------------------------------------------------------------------------
# 3. Header row (row 1) - new column names
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "ScenarioID"
$ws1.Range("B1").Value = "ChannelNameTV"
$ws1.Range("C1").Value = "TVbuyingAudience"
$ws1.Range("D1").Value = "TVSecondLengthFormat"
$ws1.Range("E1").Value = "DefaultCPMforTV"
$ws1.Range("F1").Value = "DefaultGRPs"
$ws1.Range("G1").Value = "DefaultReach"
$ws1.Range("H1").Value = "DefaultMaximumReach"
$ws1.Range("I1").Value = "DefaultPrecision"

# ---------------------------------------------------------------------------
# 4. Data row (row 2) - the single scenario-1/TV row.
#    Columns whose text looks like a plain number (ScenarioID, CPM, GRPs,
#    Reach, MaximumReach, Precision) are entered with a leading apostrophe
#    so Excel stores them as text (quoted numbers) instead of numeric
#    values, matching the original template's text-based lookup columns.
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "'1"
$ws1.Range("B2").Value = "TV"
$ws1.Range("C2").Value = "TVDONOTDELETE"
$ws1.Range("D2").Value = "TVDONOTDELETE"
$ws1.Range("E2").Value = "'20"
$ws1.Range("F2").Value = "'60"
$ws1.Range("G2").Value = "'50"
$ws1.Range("H2").Value = "'80"
$ws1.Range("I2").Value = "'70"

# ---------------------------------------------------------------------------
# 5. Cosmetic touch-ups on channel_management to line up with the new,
#    wider table: row height and column widths.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(2).RowHeight = 12.75

$ws1.Columns.Item(1).ColumnWidth = 8.74
$ws1.Columns.Item(2).ColumnWidth = 14.31
$ws1.Columns.Item(3).ColumnWidth = 14.45
$ws1.Columns.Item(4).ColumnWidth = 21.02
$ws1.Columns.Item(5).ColumnWidth = 14.74
$ws1.Columns.Item(6).ColumnWidth = 11.02
$ws1.Columns.Item(7).ColumnWidth = 11.02
$ws1.Columns.Item(8).ColumnWidth = 18.45
$ws1.Columns.Item(9).ColumnWidth = 13.88

# ---------------------------------------------------------------------------
# 6. comparison_management - rename column 1 header from "Scenario" to
#    "ScenarioID" (this also renames the Table22 list column). Nothing
#    else on this sheet changes.
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "ScenarioID"

# ---------------------------------------------------------------------------
# 7. Update the remembered cell selection on each sheet. comparison_management
#    ends up with A8 selected, channel_management ends up with J9 selected
#    and remains the active/selected tab.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A8").Select()

$ws1.Activate()
$ws1.Range("J9").Select()
